# Linha do tempo cartola 1 - add Rodada 33/34/35 scores (cols AH:AJ) and
# re-sort the standings (rows 2:21) by Total, descending - matching the
# sheet's existing AutoFilter sortState (sort key = column AN, Total).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# --- Step 1: snapshot the rows whose relative order changes once the new
# rounds are factored into the Total, BEFORE any cell is overwritten.
# (columns A..AG = team name + Rodada 01..32)
$snapshot = @{}
foreach ($r in 6,7,8,10,11,14,15) {
    $rowData = @()
    for ($c = 1; $c -le 33; $c++) {
        $rowData += $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowData
}

# --- Step 2: write each snapshot back to its new (sorted) position.
# 6<-8, 7<-6, 8<-7 (3-way rotation); 10<-11, 11<-10 and 14<-15, 15<-14 (swaps)
$moves = @{ 6 = 8; 7 = 6; 8 = 7; 10 = 11; 11 = 10; 14 = 15; 15 = 14 }
foreach ($dest in $moves.Keys) {
    $src = $moves[$dest]
    $rowData = $snapshot[$src]
    for ($c = 1; $c -le 33; $c++) {
        $ws.Cells.Item($dest, $c).Value = $rowData[$c - 1]
    }
}

# --- Step 3: enter the new Rodada 33 (AH), 34 (AI) and 35 (AJ) scores for
# every standings row, now that rows sit in their final (sorted) order.
$newScores = @{
    2  = @(61.17, 104.17, 123.61)
    3  = @(66.12, 108.71, 123.69)
    4  = @(72.57, 103.27, 123.61)
    5  = @(81.72, 81.21, 106.31)
    6  = @(79.37, 104.41, 115.14)
    7  = @(64.64, 77.01, 132.61)
    8  = @(68.52, 76.47, 119.84)
    9  = @(54.54, 56.95, 110.76)
    10 = @(73.02, 99.27, 125.44)
    11 = @(61.57, 74.72, 89.99)
    12 = @(61.97, 92.31, 109.31)
    13 = @(67.67, 52.97, 119.09)
    14 = @(65.13, 53.85, 103.69)
    15 = @(72.39, 39.80, 99.21)
    16 = @(74.67, 24.47, 67.94)
    17 = @(58.37, 47.90, 110.55)
    18 = @(36.02, 44.82, 66.39)
    19 = @(41.27, 41.50, 70.29)
    20 = @(29.67, 42.14, 50.93)
    21 = @(42.31, 30.97, 51.18)
}
foreach ($r in $newScores.Keys) {
    $vals = $newScores[$r]
    $ws.Cells.Item($r, 34).Value = $vals[0]   # AH - Rodada 33
    $ws.Cells.Item($r, 35).Value = $vals[1]   # AI - Rodada 34
    $ws.Cells.Item($r, 36).Value = $vals[2]   # AJ - Rodada 35
}

# --- Step 4: leave the selection on the Total column, like after sorting it.
[void]$ws.Range("AN2:AN21").Select()
